$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous block of 14 rows (842-855) had its "Fecha" timestamp re-saved with a
# slightly different floating point value (re-calculated on update).
for ($row = 842; $row -le 855; $row++) {
    $ws.Cells.Item($row, 4).Value = 44232.93881568287
}

# New block of 14 rows (856-869) appended after the existing data, following the
# same repeating 14-row pattern used throughout the sheet (Nombre / URL / Disponible / Fecha).
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")

# Address = target used for the relationship (without any #fragment)
# SubAddress = the fragment/location portion (only the MapStore link has one)
$addresses = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$subAddresses = @("","","","","","","","","/","","","","","")

$timestamp = 44232.9602579231
$startRow = 856

for ($i = 0; $i -lt 14; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value = $names[$i]

    $displayText = $addresses[$i]
    if ($subAddresses[$i] -ne "") {
        $displayText = $displayText + "#" + $subAddresses[$i]
    }
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.Value = $displayText
    if ($subAddresses[$i] -ne "") {
        $ws.Hyperlinks.Add($cellB, $addresses[$i], $subAddresses[$i])
    } else {
        $ws.Hyperlinks.Add($cellB, $addresses[$i])
    }
    $cellB.Style = "Hyperlink"

    $ws.Cells.Item($row, 3).Value = "Disponible"

    $cellD = $ws.Cells.Item($row, 4)
    $cellD.Value = $timestamp
    $cellD.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
